$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "CCEP" (Coca-Cola Europacific Partners) row was actually meant to be
# "KO" (The Coca-Cola Company) -- the sustainability report link belongs to
# Coca-Cola Company, not the Europacific Partners entity. Fix row 6 data.
$ws.Range("A6").Value = "KO"
$ws.Range("B6").Value = "Coca-Cola Company (The) Common Stock"
$ws.Range("C6").Value = 58.6
$ws.Range("D6").Value = -0.44
$ws.Range("E6").Value = -0.0075
$ws.Range("F6").Value = "'2.53352E+11"
$ws.Range("G6").Value = "United States"
$ws.Range("I6").Value = 65172256

# Update selection / scroll position to reflect where the user ended up
# after editing (was A8 with topLeftCell A2; now G18, no special scroll).
$ws.Range("G18").Select()
